# Reconfigure ad works dw: prefix the landing/raw/control ad_works database
# names (and the depends_on / sql references that embed them) with "yetl_".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count
$startRow = $usedRange.Row
$startCol = $usedRange.Column

for ($r = $startRow; $r -lt ($startRow + $rowCount); $r++) {
    for ($c = $startCol; $c -lt ($startCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.Contains("ad_works")) {
            $newVal = $val
            $newVal = $newVal.Replace("landing_ad_works", "yetl_landing_ad_works")
            $newVal = $newVal.Replace("raw_ad_works", "yetl_raw_ad_works")
            $newVal = $newVal.Replace("control_ad_works", "yetl_control_ad_works")
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}

# Update the window width (doubled) and move the active selection down to
# the first raw-layer row (C72), matching the author's final view state.
$excel.ActiveWindow.Width = $excel.ActiveWindow.Width * 2
$ws.Range("C72").Select()
